# "Updated for week 7" - add PCA/MDS lesson details to the schedule and
# fix a pre-existing typo ("Bonus taks" -> "Bonus task").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo on the existing "GAM and LOESS" row: "Bonus taks" -> "Bonus task"
$ws.Range("I32").Value = "Bonus task"

# --- Week 8 / Day 3 (row 40): PCA lesson - fill in notes/slides/task info.
# Note: order of assignment matters for shared-string ordering, so L/M are
# set before J to reproduce the original authoring order.
$ws.Range("F40").Value = "Notes"
$ws.Range("H40").Value = "Slides"
$ws.Range("I40").Value = "Bonus task"
$ws.Range("L40").Value = "slides/19-PCA.html"
$ws.Range("M40").Value = "#task-b8"

# --- Week 8 / Day 5 (row 43): MDS lesson - fill in notes/task info.
$ws.Range("F43").Value = "Notes"
$ws.Range("M43").Value = "#task-12"
$ws.Range("I43").Value = "Task 12"
$ws.Range("J43").Value = "#mds"

# Notes-anchor link for the PCA row, added last.
$ws.Range("J40").Value = "#PCA"

# Update the active selection to reflect where editing ended up.
$ws.Range("F43").Select()
